$d = $word.ActiveDocument

$replacements = @(
    @("2024-08-23 Friday", "2024-08-24 Saturday"),
    @("601×2=", "797×7="),
    @("381×8=", "822×6="),
    @("464×3=", "654×5="),
    @("447×8=", "285×9="),
    @("296×7=", "544×2="),
    @("825×4=", "312×2="),
    @("407×3=", "349×8="),
    @("452×2=", "710×2="),
    @("294×5=", "379×3="),
    @("199×6=", "683×7="),
    @("589×5=", "148×7="),
    @("122×5=", "646×6="),
    @("162×8=", "976×9="),
    @("299×6=", "787×4="),
    @("781×6=", "708×4="),
    @("962×5=", "142×9="),
    @("952×3=", "246×6="),
    @("795×2=", "247×7="),
    @("970×7=", "283×3="),
    @("974×8=", "494×9="),
    @("669×6=", "678×4="),
    @("588×5=", "241×3="),
    @("595×5=", "591×4="),
    @("109×7=", "681×8="),
    @("864×7=", "696×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
